# Weekly update: insert a new observation row at row 83 (Fecha 2021-09-27)
# for "Hortaliza, Femacal de La Calera - Achicoria", shifting the existing
# rows 83-126 down to 84-127.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 83, pushing rows 83..126 to 84..127.
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with the new record's data.
$ws.Range("A83").Value = 3
$ws.Range("B83").Value = "Femacal de La Calera"
$ws.Range("C83").Value = "Coquimbo"
$ws.Range("D83").Value = 44466
$ws.Range("E83").Value = 5
$ws.Range("F83").Value = 100112010
$ws.Range("G83").Value = "Achicoria"
$ws.Range("H83").Value = "Sin especificar"
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value = 130
$ws.Range("K83").Value = 5500
$ws.Range("L83").Value = 6000
$ws.Range("M83").Value = 5731
$ws.Range("N83").Value = "`$/caja 16 unidades"
$ws.Range("O83").Value = "Provincia de Quillota"
$ws.Range("P83").Value = 358
$ws.Range("Q83").Value = 16
$ws.Range("R83").Value = "Hortaliza"
